# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Periodo Mora" column (E16:E20): previously 1709,1712,1801,1802,1803
# now listed 1803,1802,1801,1712,1709 (newest first).
$ws.Range("E16").Value = "1803"
$ws.Range("E17").Value = "1802"
$ws.Range("E18").Value = "1801"
$ws.Range("E19").Value = "1712"
$ws.Range("E20").Value = "1709"

# Valor Mora column (F16:F20): the 19673 figure that used to belong to the
# last period (1803, row 20) now belongs to the first period (1803, row 16);
# the other rows keep 29509.
$ws.Range("F16").Value = 19673
$ws.Range("F17").Value = 29509
$ws.Range("F18").Value = 29509
$ws.Range("F19").Value = 29509
$ws.Range("F20").Value = 29509
